$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = [double]"2.900324425070266e-11"
$ws.Range("C2").Value = [double]"1.27202479660582e-08"
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 9.564027516723449

# Row 3
$ws.Range("B3").Value = [double]"1.063418937352623e-07"
$ws.Range("C3").Value = 0.0001537489499301437
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 9.564181359266023
